{"js": "const replacements = [\n  [\"53\u00f72=\", \"28\u00f74=\"],\n  [\"96\u00f77=\", \"84\u00f78=\"],\n  [\"22\u00f73=\", \"43\u00f79=\"],\n  [\"32\u00f78=\", \"92\u00f73=\"],\n  [\"77\u00f75=\", \"84\u00f75=\"],\n  [\"41\u00f78=\", \"87\u00f78=\"],\n  [\"49\u00f77=\", \"28\u00f76=\"],\n  [\"37\u00f75=\", \"75\u00f73=\"],\n  [\"41\u00f76=\", \"89\u00f73=\"],\n  [\"71\u00f78=\", \"54\u00f79=\"],\n  [\"56\u00f78=\", \"60\u00f72=\"],\n  [\"95\u00f75=\", \"99\u00f78=\"],\n  [\"47\u00f75=\", \"49\u00f72=\"],\n  [\"22\u00f72=\", \"35\u00f74=\"],\n  [\"11\u00f74=\", \"61\u00f79=\"],\n  [\"23\u00f73=\", \"37\u00f72=\"],\n  [\"43\u00f76=\", \"17\u00f72=\"],\n  [\"89\u00f78=\", \"63\u00f74=\"],\n  [\"41\u00f77=\", \"30\u00f79=\"],\n  [\"26\u00f78=\", \"81\u00f75=\"],\n  [\"85\u00f75=\", \"57\u00f76=\"],\n  [\"21\u00f75=\", \"92\u00f73=\"],\n  [\"32\u00f73=\", \"75\u00f78=\"],\n  [\"39\u00f73=\", \"92\u00f73=\"],\n  [\"16\u00f74=\", \"73\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"53\u00f72=\", \"28\u00f74=\"),\n    @(\"96\u00f77=\", \"84\u00f78=\"),\n    @(\"22\u00f73=\", \"43\u00f79=\"),\n    @(\"32\u00f78=\", \"92\u00f73=\"),\n    @(\"77\u00f75=\", \"84\u00f75=\"),\n    @(\"41\u00f78=\", \"87\u00f78=\"),\n    @(\"49\u00f77=\", \"28\u00f76=\"),\n    @(\"37\u00f75=\", \"75\u00f73=\"),\n    @(\"41\u00f76=\", \"89\u00f73=\"),\n    @(\"71\u00f78=\", \"54\u00f79=\"),\n    @(\"56\u00f78=\", \"60\u00f72=\"),\n    @(\"95\u00f75=\", \"99\u00f78=\"),\n    @(\"47\u00f75=\", \"49\u00f72=\"),\n    @(\"22\u00f72=\", \"35\u00f74=\"),\n    @(\"11\u00f74=\", \"61\u00f79=\"),\n    @(\"23\u00f73=\", \"37\u00f72=\"),\n    @(\"43\u00f76=\", \"17\u00f72=\"),\n    @(\"89\u00f78=\", \"63\u00f74=\"),\n    @(\"41\u00f77=\", \"30\u00f79=\"),\n    @(\"26\u00f78=\", \"81\u00f75=\"),\n    @(\"85\u00f75=\", \"57\u00f76=\"),\n    @(\"21\u00f75=\", \"92\u00f73=\"),\n    @(\"32\u00f73=\", \"75\u00f78=\"),\n    @(\"39\u00f73=\", \"92\u00f73=\"),\n    @(\"16\u00f74=\", \"73\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
